$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E4").Value = 52
$ws.Range("E5").Value = 163
$ws.Range("F5").Value = 111
$ws.Range("H5").Value = 122
$ws.Range("E7").Value = 44
$ws.Range("F7").Value = 30
$ws.Range("H7").Value = 34
$ws.Range("E8").Value = 10
$ws.Range("E10").Value = 700
$ws.Range("F10").Value = 390
$ws.Range("H10").Value = 485
$ws.Range("E11").Value = 461
$ws.Range("F11").Value = 261
$ws.Range("H11").Value = 326
$ws.Range("E12").Value = 699
$ws.Range("F12").Value = 422
$ws.Range("H12").Value = 508
$ws.Range("E13").Value = 166
$ws.Range("E14").Value = 143
$ws.Range("F14").Value = 83
$ws.Range("H14").Value = 117
$ws.Range("E15").Value = 203
$ws.Range("E18").Value = 63
$ws.Range("F18").Value = 35
$ws.Range("H18").Value = 52
$ws.Range("F20").Value = 46
$ws.Range("H20").Value = 83
$ws.Range("E22").Value = 196
$ws.Range("E23").Value = 227
$ws.Range("E24").Value = 271
$ws.Range("F24").Value = 158
$ws.Range("H24").Value = 188
$ws.Range("E25").Value = 331
$ws.Range("F25").Value = 183
$ws.Range("H25").Value = 243
$ws.Range("E26").Value = 198
$ws.Range("E27").Value = 382
$ws.Range("F27").Value = 210
$ws.Range("H27").Value = 292
$ws.Range("E28").Value = 228
$ws.Range("E29").Value = 198
$ws.Range("E30").Value = 253
$ws.Range("F30").Value = 158
$ws.Range("H30").Value = 210
$ws.Range("E31").Value = 84
$ws.Range("F31").Value = 38
$ws.Range("H31").Value = 66
$ws.Range("E32").Value = 217
$ws.Range("F32").Value = 138
$ws.Range("H32").Value = 176
$ws.Range("E33").Value = 328
$ws.Range("E34").Value = 251
$ws.Range("E37").Value = 197
$ws.Range("F37").Value = 111
$ws.Range("H37").Value = 148
$ws.Range("E39").Value = 200
$ws.Range("E40").Value = 308
$ws.Range("F40").Value = 157
$ws.Range("H40").Value = 237
$ws.Range("E41").Value = 443
$ws.Range("F41").Value = 222
$ws.Range("H41").Value = 314
$ws.Range("E42").Value = 467
$ws.Range("F42").Value = 265
$ws.Range("H42").Value = 326
$ws.Range("E43").Value = 143
$ws.Range("F43").Value = 81
$ws.Range("H43").Value = 108
$ws.Range("E44").Value = 371
$ws.Range("F44").Value = 194
$ws.Range("H44").Value = 262
$ws.Range("E45").Value = 181
$ws.Range("E46").Value = 391
$ws.Range("F46").Value = 229
$ws.Range("H46").Value = 293
$ws.Range("E47").Value = 539
$ws.Range("F47").Value = 305
$ws.Range("H47").Value = 397
$ws.Range("E48").Value = 270
$ws.Range("F48").Value = 131
$ws.Range("H48").Value = 175
$ws.Range("E49").Value = 339
$ws.Range("E50").Value = 284
$ws.Range("F50").Value = 157
$ws.Range("H50").Value = 230
